# Saldo.xlsx update:
#  - Remove rows for accounts 005660459 (Moacir), 004503381 (Frederico),
#    005886225 (Vinicius), 004216504 (Wander), 004426743 (Gabrielle)
#  - Change the account number on Andrea's row from 003894173 to 005186167

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Account numbers (column A, inline strings) whose rows must be deleted entirely.
$accountsToDelete = @("005660459", "004503381", "005886225", "004216504", "004426743")

# Resolve each account number to its current row, then delete highest row
# number first so earlier matches' row numbers stay valid as rows shift up.
$rowsToDelete = @()
foreach ($acct in $accountsToDelete) {
    $cell = $ws.Cells.Find($acct)
    if ($cell -ne $null) {
        $rowsToDelete += $cell.Row
    }
}

$rowsToDelete = $rowsToDelete | Sort-Object -Descending

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Update Andrea's account number. Force text format first so the
# leading zeros in the new account number are preserved.
$andreaCell = $ws.Cells.Find("003894173")
if ($andreaCell -ne $null) {
    $andreaCell.NumberFormat = "@"
    $andreaCell.Value = "005186167"
}
